$d = $word.ActiveDocument

# The document contains five "<id>...</id>" tags (for p111v_1 .. p111v_5),
# each currently split across three separate runs:
#   run1: "<id>"     (Courier New, color 7f6000, sz 18 - "tag" formatting)
#   run2: "p111v_N"  (plain color 000000 - "value" formatting)
#   run3: "</id>"    (Courier New, color 7f6000, sz 18 - "tag" formatting)
#
# Re-downloading collapses each triple into a single run (keeping the
# formatting of the first/opening-tag run) containing the full
# "<id>p111v_N</id>" text. Replacing the found text with itself via
# Find.Execute merges the matched runs into one run using the formatting
# of the first run in the match, exactly reproducing that collapse.

for ($n = 1; $n -le 5; $n++) {
    $needle = "<id>p111v_$n</id>"
    $found = $d.Content.Find.Execute($needle, $true, $true, $false, $false, $false, $true, 1, $false, $needle, 2)
    if (-not $found) {
        Write-Host "WARNING: pattern not found for p111v_$n"
    }
}
